$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all existing data rows (2-523)
$ws.Range("C2:C523").Value = 45175

# Row 523 picks up an explicit row height (matches the newly appended row below it)
$ws.Rows.Item(523).RowHeight = 15

# Append the new record as row 524
$ws.Range("A524").Value = "A 41065-2023"

$ws.Range("B524").NumberFormat = "YYYY-MM-DD"
$ws.Range("B524").Value = 45173

$ws.Range("C524").NumberFormat = "YYYY-MM-DD"
$ws.Range("C524").Value = 45175

$ws.Range("D524").Value = "SKÅNE LÄN"
$ws.Range("E524").Value = "KRISTIANSTAD"

$ws.Range("G524").Value = 1
$ws.Range("H524").Value = 0
$ws.Range("I524").Value = 0
$ws.Range("J524").Value = 0
$ws.Range("K524").Value = 0
$ws.Range("L524").Value = 0
$ws.Range("M524").Value = 0
$ws.Range("N524").Value = 0
$ws.Range("O524").Value = 0
$ws.Range("P524").Value = 0
$ws.Range("Q524").Value = 0

$ws.Range("R524").WrapText = $true
$ws.Range("R524").Value = ""
